$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C5").Value = 72607.85767434421
$ws.Range("C6").Value = 72607.85767434421
$ws.Range("C7").Value = 70429.62194411388
$ws.Range("C8").Value = 12905.604041271381
$ws.Range("C9").Value = 20208.714613346274
$ws.Range("C10").Value = 52399.14306099794
$ws.Range("C11").Value = 52399.14306099794
$ws.Range("C14").Value = 38749.14306099795
$ws.Range("C15").Value = 38003.55461973628
$ws.Range("C16").Value = 363.0156987616624
$ws.Range("C19").Value = 36420.47031849796
$ws.Range("C20").Value = 21733.049112279747

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 7260.313975233249
$ws.Range("C3").Value = 8320.666666666664
$ws.Range("D3").Value = 14.604777356055736
$ws.Range("C5").Value = 8320.666666666664
$ws.Range("C8").Value = 10957.0
$ws.Range("D8").Value = 50.916338293041804
$ws.Range("C9").Value = 6396.0
$ws.Range("D9").Value = -11.904636330903044
$ws.Range("C10").Value = 16941.0
$ws.Range("D10").Value = 133.33701624736892
$ws.Range("C11").Value = 6894.0
$ws.Range("D11").Value = -5.045428840720071
$ws.Range("C12").Value = 6590.0
$ws.Range("D12").Value = -9.232575581715299
$ws.Range("C13").Value = 11196.0
$ws.Range("D13").Value = 54.20820694796897
$ws.Range("C14").Value = 22194.0
$ws.Range("D14").Value = 205.68925911068447
$ws.Range("D15").Value = 8.686759648662305

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 7695.932813747244
$ws.Range("C3").Value = 7260.25
$ws.Range("D3").Value = -5.66120864476602
$ws.Range("C5").Value = 7260.249999999999
$ws.Range("A8").Value = "TORENBEEK_1982"
$ws.Range("C8").Value = 6357.0
$ws.Range("D8").Value = -17.397927530701782
$ws.Range("A9").Value = "SADRAEY"
$ws.Range("C9").Value = 8723.0
$ws.Range("D9").Value = 13.34558410408815
$ws.Range("A10").Value = "ROSKAM"
$ws.Range("C10").Value = 7694.0
$ws.Range("D10").Value = -0.025114742995044625
$ws.Range("C11").Value = 6267.0
$ws.Range("D11").Value = -18.567376409455413

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 834.9361071518235
$ws.Range("C3").Value = 775.75
$ws.Range("D3").Value = -7.088698960896815
$ws.Range("C5").Value = 775.7499999999999
$ws.Range("A8").Value = "TORENBEEK_1976"
$ws.Range("C8").Value = 1040.0
$ws.Range("D8").Value = 24.56042936599075
$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 1640.0
$ws.Range("D9").Value = 96.42221553867772
$ws.Range("A10").Value = "NICOLAI_2013"
$ws.Range("C10").Value = 439.0
$ws.Range("D10").Value = -47.42112645031737
$ws.Range("A11").Value = "SADRAEY"
$ws.Range("C11").Value = 551.0
$ws.Range("D11").Value = -34.00692636474913
$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 740.0
$ws.Range("D12").Value = -11.370463720352735
$ws.Range("A13").Value = "NICOLAI_1984"
$ws.Range("C13").Value = -425.0
$ws.Range("D13").Value = -150.9020985389866
$ws.Range("A14").Value = "HOWE"
$ws.Range("C14").Value = 1521.0
$ws.Range("D14").Value = 82.16962794776147
$ws.Range("A15").Value = "RAYMER"
$ws.Range("C15").Value = 700.0
$ws.Range("D15").Value = -16.161249465198534

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 834.9361071518235
$ws.Range("C3").Value = 661.5714285714284
$ws.Range("D3").Value = -20.763825770068255
$ws.Range("C5").Value = 661.5714285714284
$ws.Range("A8").Value = "TORENBEEK_1976"
$ws.Range("C8").Value = 749.0
$ws.Range("D8").Value = -10.292536927762429
$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 1640.0
$ws.Range("D9").Value = 96.42221553867772
$ws.Range("A10").Value = "SADRAEY"
$ws.Range("C10").Value = 194.0
$ws.Range("D10").Value = -76.76468913749788
$ws.Range("A11").Value = "ROSKAM"
$ws.Range("C11").Value = 511.0
$ws.Range("D11").Value = -38.797712109594926
$ws.Range("A12").Value = "NICOLAI_1984"
$ws.Range("C12").Value = -195.0
$ws.Range("D12").Value = -123.35508050612326
$ws.Range("A13").Value = "HOWE"
$ws.Range("C13").Value = 1230.0
$ws.Range("D13").Value = 47.31666165400829
$ws.Range("A14").Value = "RAYMER"
$ws.Range("C14").Value = 502.0
$ws.Range("D14").Value = -39.87563890218523

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1379.4596552943171
$ws.Range("D3").Value = 101.67316886164687
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 1374.0
$ws.Range("D10").Value = 99.20843566923278
$ws.Range("A11").Value = "KUNDU"
$ws.Range("C11").Value = 1389.0
$ws.Range("D11").Value = 101.38320025077462
$ws.Range("A12").Value = "RAYMER"
$ws.Range("C12").Value = 1410.0
$ws.Range("D12").Value = 104.42787066493321
$ws.Range("A17").Value = "KROO"
$ws.Range("C17").Value = 1374.0
$ws.Range("D17").Value = 99.20843566923278
$ws.Range("A18").Value = "KUNDU"
$ws.Range("C18").Value = 1389.0
$ws.Range("D18").Value = 101.38320025077462
$ws.Range("A19").Value = "RAYMER"
$ws.Range("C19").Value = 1410.0
$ws.Range("D19").Value = 104.42787066493321

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 6026.060599443597
$ws.Range("D3").Value = 7.057226306834707
$ws.Range("A11").Value = "NICOLAI_1984"
$ws.Range("D11").Value = -1.9591671456887836
$ws.Range("D12").Value = 8.362667322046757
$ws.Range("D13").Value = 14.768178744146304
$ws.Range("A18").Value = "NICOLAI_1984"
$ws.Range("D18").Value = -1.9591671456887836
$ws.Range("D19").Value = 8.362667322046757
$ws.Range("D20").Value = 14.768178744146304

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2976.728729845632
$ws.Range("C3").Value = 1932.8110170416635
$ws.Range("D3").Value = -35.069292755410196
$ws.Range("C5").Value = 1932.8110170416633
$ws.Range("A9").Value = "NICOLAI_1984"
$ws.Range("C9").Value = 1932.8110170416635
$ws.Range("D9").Value = -35.06929275541022
$ws.Range("A11").Value = "NICOLAI_1984"
$ws.Range("C11").Value = 226.12722719328167
$ws.Range("A13").Value = "NICOLAI_1984"
$ws.Range("C13").Value = 1706.683789848382

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 9874.027006317217
$ws.Range("C3").Value = 8236.087872884898
$ws.Range("D3").Value = -16.588359869629663
$ws.Range("C4").Value = 8236.087872884897
$ws.Range("C8").Value = 8236.087872884898
$ws.Range("D8").Value = -16.58835986962965
$ws.Range("C21").Value = 1049.4473666076028
$ws.Range("C23").Value = 1049.4473666076026
$ws.Range("C26").Value = 546.2874086241745
$ws.Range("C28").Value = 546.2874086241744
$ws.Range("C36").Value = 855.7018854509024
$ws.Range("C38").Value = 855.7018854509023
$ws.Range("C41").Value = 3213.3940371457693
$ws.Range("C43").Value = 3213.3940371457693
